$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.736.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.599.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.619.12"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.369"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.057.61"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.697.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.611.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.58"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.78%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.519"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +15.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.54"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.19"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.916"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.76"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.76"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.33%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.601"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.84%  "
